# Add a newly completed book to the "Completed" sheet's reading log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# The sheet currently has data through row 74 (row 1 is the header row).
# Append the new book as row 75.
$newRow = 75

$ws.Cells.Item($newRow, 1).Value = "From Bacteria to Bach and Back"
$ws.Cells.Item($newRow, 2).Value = "Daniel C Dennett"
$ws.Cells.Item($newRow, 3).Value = 43959
$ws.Cells.Item($newRow, 4).Value = 43962
$ws.Cells.Item($newRow, 5).Value = "consciousness;matter;science;neurology;free will"
$ws.Cells.Item($newRow, 6).Value = "Audio"
$ws.Cells.Item($newRow, 7).Value = "15 Hours 52 Mins"

# Match date formatting used by the existing Start Date / Finish Date columns
# by copying the format already applied to the row above, rather than
# introducing a brand new number format.
$ws.Range("C74:D74").Copy()
$ws.Range("C75:D75").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Reflect the scrolled/selected state left behind after entering the row.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 56
$ws.Range("E76").Select()
